$wb = $excel.ActiveWorkbook

$sprint = $wb.Worksheets.Item("Sprint")

# Update F8: task status from "In progress" to "Done"
$sprint.Range("F8").Value = "Done"

# Set R8 = 2 (effort logged on day 9)
$sprint.Range("R8").Value = 2

# Activate Sprint sheet, select R9 (mirrors the user's click back onto
# the Sprint tab after logging today's effort)
$sprint.Activate()
$null = $sprint.Range("R9").Select()

